# Update Name of Algo
# Apply updated values to column B per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value  = 6.152499999999999
$ws.Range("B18").Value = 6.477399999999998
$ws.Range("B20").Value = 8.910199999999996
$ws.Range("B27").Value = 6.420600000000007
$ws.Range("B35").Value = 8.6991
$ws.Range("B69").Value = 5.401399999999994
$ws.Range("B76").Value = 5.477199999999999
$ws.Range("B78").Value = 10.2266
$ws.Range("B82").Value = 5.417500000000001
$ws.Range("B83").Value = 5.300699999999997
$ws.Range("B93").Value = 5.599100000000001
